$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4400
$ws.Range("J64").Value = 4400
$ws.Range("L64").Value = 4400
$ws.Range("N64").Value = -4896

$ws.Range("H67").Value = 4400
$ws.Range("J67").Value = 4400
$ws.Range("L67").Value = 4400
$ws.Range("N67").Value = -6116

$ws.Range("H70").Value = 1574.7059
$ws.Range("I70").Value = 1472
$ws.Range("J70").Value = 1721.4286
$ws.Range("K70").Value = 4416
$ws.Range("L70").Value = 5164.2858
$ws.Range("M70").Value = -4146
$ws.Range("N70").Value = -5704.2858

$ws.Range("H73").Value = 1574.7059
$ws.Range("I73").Value = 1472
$ws.Range("J73").Value = 1721.4286
$ws.Range("K73").Value = 4416
$ws.Range("L73").Value = 5164.2858
$ws.Range("M73").Value = -3480
$ws.Range("N73").Value = -7036.2858

$ws.Range("H92").Value = 1168
$ws.Range("I92").Value = 1168
$ws.Range("K92").Value = 1168
$ws.Range("M92").Value = 80

$ws.Range("H115").Value = 3399.6667
$ws.Range("I115").Value = 2999
$ws.Range("K115").Value = 8997
$ws.Range("M115").Value = -7430

$ws.Range("H137").Value = 2719.4
$ws.Range("I137").Value = 2719.4
$ws.Range("K137").Value = 8158.200000000001
$ws.Range("M137").Value = -5608.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2222
$ws.Range("I74").Value = 2222
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2222
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1348
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 2222
$ws.Range("I77").Value = 2222
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 11110
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -6742
$ws.Range("N77").ClearContents()

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H130").Value = 24145
$ws.Range("J130").Value = 24145
$ws.Range("L130").Value = 24145
$ws.Range("N130").Value = -34185

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 7375
$ws.Range("I38").Value = 7000
$ws.Range("K38").Value = 7000
$ws.Range("M38").Value = -6584

$ws.Range("H80").Value = 1246
$ws.Range("I80").Value = 776.3333
$ws.Range("K80").Value = 776.3333
$ws.Range("M80").Value = 221.6667

$ws.Range("H83").Value = 1246
$ws.Range("I83").Value = 776.3333
$ws.Range("K83").Value = 3881.6665
$ws.Range("M83").Value = 1110.3335

$ws.Range("H94").Value = 4502.5
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 4502.5
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 4502.5
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -5404.5

$ws.Range("H99").Value = 1048.3334
$ws.Range("I99").Value = 897.5
$ws.Range("J99").Value = 1123.75
$ws.Range("K99").Value = 897.5
$ws.Range("L99").Value = 1123.75
$ws.Range("M99").Value = 600.5
$ws.Range("N99").Value = -4119.75

$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 533.3333
$ws.Range("I22").Value = 450
$ws.Range("J22").Value = 700
$ws.Range("K22").Value = 450
$ws.Range("L22").Value = 700
$ws.Range("M22").Value = -100
$ws.Range("N22").Value = -1400

$ws.Range("H62").Value = 9166.5
$ws.Range("J62").Value = 13333
$ws.Range("L62").Value = 13333
$ws.Range("N62").Value = -14581

$ws.Range("H65").Value = 9166.5
$ws.Range("J65").Value = 13333
$ws.Range("L65").Value = 66665
$ws.Range("N65").Value = -72905

$ws.Range("H86").Value = 3334.3333
$ws.Range("I86").Value = 3497.5
$ws.Range("K86").Value = 3497.5
$ws.Range("M86").Value = -2374.5

$ws.Range("H89").Value = 3334.3333
$ws.Range("I89").Value = 3497.5
$ws.Range("K89").Value = 17487.5
$ws.Range("M89").Value = -11871.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 165
$ws.Range("N13").ClearContents()

$ws.Range("H18").Value = 946
$ws.Range("I18").Value = 946
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 2838
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -2669
$ws.Range("N18").ClearContents()

$ws.Range("H46").Value = 200
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H109").Value = 2273.3
$ws.Range("I109").Value = 1303.6666
$ws.Range("K109").Value = 3910.9998
$ws.Range("M109").Value = -2870.9998

$ws.Range("H132").Value = 1136
$ws.Range("I132").Value = 1136
$ws.Range("K132").Value = 10224
$ws.Range("M132").Value = -7694

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 4556097
$ws.Range("I11").Value = 4250719
$ws.Range("K11").Value = 4250719
$ws.Range("M11").Value = -4250580

$ws.Range("H102").Value = 433
$ws.Range("I102").Value = 433
$ws.Range("K102").Value = 433
$ws.Range("M102").Value = 1189

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 75895.5
$ws.Range("J136").Value = 121655.836
$ws.Range("L136").Value = 364967.508
$ws.Range("N136").Value = -370067.508

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3995
$ws.Range("I96").Value = 4500
$ws.Range("K96").Value = 4500
$ws.Range("M96").Value = -3127

$ws.Range("H105").Value = 58975
$ws.Range("J105").Value = 58975
$ws.Range("L105").Value = 58975
$ws.Range("N105").Value = -65963

$ws.Range("H107").Value = 1421.1
$ws.Range("I107").Value = 870.8
$ws.Range("K107").Value = 2612.4
$ws.Range("M107").Value = -692.3999999999996

$ws.Range("H132").Value = 1586.875
$ws.Range("I132").Value = 1242.2858
$ws.Range("K132").Value = 3926.8574
$ws.Range("M132").Value = -1196.8574
